$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.093.29'
$ws.Range('E2').Value = '  +1.64%  '
$ws.Range('D3').Value = '2.018.38'
$ws.Range('E3').Value = '  +3.61%  '
$ws.Range('E4').Value = '  -0.06%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '246.78'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.54%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.627'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +1.10%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '60.03'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +0.48%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.391'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +3.73%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0814'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +3.59%  '
$ws.Range('E11').Value = '  +1.90%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '15.10'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +6.96%  '
$ws.Range('D13').Value = '2.318.50'
$ws.Range('E13').Value = '  +3.65%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.851'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +3.60%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '21.77'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +1.96%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '5.44'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +4.09%  '
$ws.Range('D17').Value = '2.020.36'
$ws.Range('E17').Value = '  +3.48%  '
$ws.Range('D18').Value = '37.088.85'
$ws.Range('E18').Value = '  +1.90%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '70.34'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +1.81%  '
$ws.Range('D20').Value = '0.0₃0867'
$ws.Range('E20').Value = '  +2.17%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '5.22'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +3.27%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '230.37'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.77%  '
$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '2.61'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +7.52%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -0.03%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.36'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.50%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '9.40'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +3.25%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '162.78'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +1.88%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '0.137'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -4.79%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '19.70'
$c.Style = 'Normal'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.37'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +5.39%  '
$ws.Range('E31').Value = '  +1.53%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '4.77'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +0.93%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.0670'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +9.83%  '
$ws.Range('E34').Value = '  +10.02%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '4.45'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.48%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '3.62'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +5.50%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('E38').Value = '  +2.04%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '5.37'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -0.28%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '3.05'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +4.60%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.0978'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +1.13%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '16.93'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +8.06%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.18'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +2.17%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.0213'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +2.34%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '91.59'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +3.61%  '
$ws.Range('D46').Value = '1.380.77'
$ws.Range('E46').Value = '  +1.52%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.05'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +3.10%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '7.41'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +4.34%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '2.14'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +17.50%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '2.87'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +1.72%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '46.02'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +2.10%  '
